$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column R (year 2021) to the header row, matching the format
# already used for the neighboring Q4 cell.
$ws.Range("Q4").Copy()
$ws.Range("R4").PasteSpecial(-4122)
$ws.Range("R4").Value = 2021

# Add the corresponding data value in row 5, matching the "0.0" number
# format used elsewhere in that data row (e.g. H5/M5).
$ws.Range("H5").Copy()
$ws.Range("R5").PasteSpecial(-4122)
$ws.Range("R5").Value = 18.953297329007047

# Update the active selection to Q8 (as recorded in the saved view state).
$ws.Activate()
$ws.Range("Q8").Select()
